$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.971.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.532.78"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.29%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.21"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.531.45"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.34%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.66%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.92"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.384"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.127.85"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.18"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000181"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.527.07"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.03%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.803.96"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.04"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.75%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.22"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.68"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.574"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.666.53"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.89"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000114"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.67"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.26"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.19"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.540.98"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.74"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.36"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +14.26%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "170.07"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.55"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.86"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.94"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0800"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.820"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.74"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +17.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.55"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.42"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.24%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.66"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.86"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.407.48"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +11.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "305.85"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +12.41%  "
